$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.927247083433903
$ws.Range("D2").Value = 3.775743729127368
$ws.Range("E2").Value = 40.45288384109088
$ws.Range("F2").Value = 14.97546744076235
$ws.Range("G2").Value = 13.46280818360098
$ws.Range("H2").Value = 10.01677876861353
$ws.Range("M2").Value = 57.83902694580253
$ws.Range("O2").Value = 13.17240368640964

$ws.Range("C3").Value = 4.750517232111898
$ws.Range("D3").Value = 3.661503447323658
$ws.Range("E3").Value = 37.73155792820933
$ws.Range("F3").Value = 15.12409468051381
$ws.Range("G3").Value = 13.55267791125082
$ws.Range("H3").Value = 10.16038209066586
$ws.Range("M3").Value = 54.25941303192177
$ws.Range("O3").Value = 13.40649720313919

$ws.Range("C4").Value = 4.640260160178539
$ws.Range("D4").Value = 3.58900865225798
$ws.Range("E4").Value = 35.96029272391751
$ws.Range("F4").Value = 15.23188147584457
$ws.Range("G4").Value = 13.64224180635605
$ws.Range("H4").Value = 10.25347861643226
$ws.Range("M4").Value = 51.93329055210783
$ws.Range("O4").Value = 13.56205856340795

$ws.Range("C5").Value = 4.594969629215675
$ws.Range("D5").Value = 3.55890075111013
$ws.Range("E5").Value = 35.21310432235119
$ws.Range("F5").Value = 15.27981277178268
$ws.Range("G5").Value = 13.68695014781933
$ws.Range("H5").Value = 10.29263916969765
$ws.Range("M5").Value = 50.95313951435505
$ws.Range("O5").Value = 13.62833172933574

$ws.Range("C6").Value = 4.587429866860584
$ws.Range("D6").Value = 3.553867989574305
$ws.Range("E6").Value = 35.0874970203289
$ws.Range("F6").Value = 15.28800942053293
$ws.Range("G6").Value = 13.69485660623713
$ws.Range("H6").Value = 10.29921519076207
$ws.Range("M6").Value = 50.7884402606981
$ws.Range("O6").Value = 13.63950770046731

$ws.Range("C7").Value = 4.639650704417154
$ws.Range("D7").Value = 3.5886048613432
$ws.Range("E7").Value = 35.95031884581027
$ws.Range("F7").Value = 15.23251186480113
$ws.Range("G7").Value = 13.64281210950797
$ws.Range("H7").Value = 10.25400181681712
$ws.Range("M7").Value = 51.9202023776633
$ws.Range("O7").Value = 13.56294080203689

$ws.Range("C8").Value = 4.866717782484991
$ws.Range("D8").Value = 3.736855987602195
$ws.Range("E8").Value = 39.53530928296156
$ws.Range("F8").Value = 15.02319363143655
$ws.Range("G8").Value = 13.48640607599351
$ws.Range("H8").Value = 10.06526194991403
$ws.Range("M8").Value = 56.63137824753908
$ws.Range("O8").Value = 13.25060652299585

$ws.Range("C9").Value = 5.295044451952883
$ws.Range("D9").Value = 4.007948417325459
$ws.Range("E9").Value = 45.77616986598547
$ws.Range("F9").Value = 14.75053714856692
$ws.Range("G9").Value = 13.47003255103211
$ws.Range("H9").Value = 9.734934612170903
$ws.Range("M9").Value = 64.85462759664961
$ws.Range("O9").Value = 12.73644523263456

$ws.Range("C10").Value = 5.59572208193655
$ws.Range("D10").Value = 4.194041492916691
$ws.Range("E10").Value = 49.8903129724787
$ws.Range("F10").Value = 14.64325077443523
$ws.Range("G10").Value = 13.65503134276623
$ws.Range("H10").Value = 9.517583551060333
$ws.Range("M10").Value = 70.28269789865521
$ws.Range("O10").Value = 12.42527035288684

$ws.Range("C11").Value = 5.728821765908718
$ws.Range("D11").Value = 4.275675891399866
$ws.Range("E11").Value = 51.66196655539371
$ws.Range("F11").Value = 14.61651868239649
$ws.Range("G11").Value = 13.78491414824393
$ws.Range("H11").Value = 9.424473103467905
$ws.Range("M11").Value = 72.62056413221651
$ws.Range("O11").Value = 12.29974885781246

$ws.Range("C12").Value = 5.778648805955622
$ws.Range("D12").Value = 4.306142030171384
$ws.Range("E12").Value = 52.31867399223526
$ws.Range("F12").Value = 14.60971724218478
$ws.Range("G12").Value = 13.84081122365048
$ws.Range("H12").Value = 9.39006827578269
$ws.Range("M12").Value = 73.48712285430619
$ws.Range("O12").Value = 12.2546647024433

$ws.Range("C13").Value = 5.767943897807394
$ws.Range("D13").Value = 4.299600663874775
$ws.Range("E13").Value = 52.17786792005437
$ws.Range("F13").Value = 14.61103206017185
$ws.Range("G13").Value = 13.82847284385426
$ws.Range("H13").Value = 9.397439559676844
$ws.Range("M13").Value = 73.30132479105309
$ws.Range("O13").Value = 12.26426314382847

$ws.Range("C14").Value = 5.73293283009596
$ws.Range("D14").Value = 4.278191383472222
$ws.Range("E14").Value = 51.71627736709515
$ws.Range("F14").Value = 14.61589173786709
$ws.Range("G14").Value = 13.7893777763891
$ws.Range("H14").Value = 9.42162532799143
$ws.Range("M14").Value = 72.69223108714959
$ws.Range("O14").Value = 12.29598984755209

$ws.Range("C15").Value = 5.711411410469029
$ws.Range("D15").Value = 4.265019022935601
$ws.Range("E15").Value = 51.43169884217007
$ws.Range("F15").Value = 14.6193052357281
$ws.Range("G15").Value = 13.76630795241545
$ws.Range("H15").Value = 9.436551823147465
$ws.Range("M15").Value = 72.31670788981579
$ws.Range("O15").Value = 12.31574655758424

$ws.Range("C16").Value = 5.586945608787239
$ws.Range("D16").Value = 4.188644579936445
$ws.Range("E16").Value = 49.77253685757127
$ws.Range("F16").Value = 14.64545469026873
$ws.Range("G16").Value = 13.64747688454793
$ws.Range("H16").Value = 9.523786628112106
$ws.Range("M16").Value = 70.12728163105348
$ws.Range("O16").Value = 12.43380810536653

$ws.Range("C17").Value = 5.509612899723999
$ws.Range("D17").Value = 4.14100801933536
$ws.Range("E17").Value = 48.72924001855102
$ws.Range("F17").Value = 14.66725681720675
$ws.Range("G17").Value = 13.58641215692899
$ws.Range("H17").Value = 9.578795867662611
$ws.Range("M17").Value = 68.75059068733704
$ws.Range("O17").Value = 12.51044036807722

$ws.Range("C18").Value = 5.464789130249715
$ws.Range("D18").Value = 4.113325064446339
$ws.Range("E18").Value = 48.11974967128414
$ws.Range("F18").Value = 14.68186597927896
$ws.Range("G18").Value = 13.55558991630652
$ws.Range("H18").Value = 9.610976589712186
$ws.Range("M18").Value = 67.9463770181391
$ws.Range("O18").Value = 12.55601732335547

$ws.Range("C19").Value = 5.44955497273042
$ws.Range("D19").Value = 4.103903795292363
$ws.Range("E19").Value = 47.91176437889001
$ws.Range("F19").Value = 14.68716317487092
$ws.Range("G19").Value = 13.54588718514399
$ws.Range("H19").Value = 9.621964587878194
$ws.Range("M19").Value = 67.67195387276168
$ws.Range("O19").Value = 12.57170233570131

$ws.Range("C20").Value = 5.517881055355305
$ws.Range("D20").Value = 4.146108469856508
$ws.Range("E20").Value = 48.8412736988415
$ws.Range("F20").Value = 14.66472070554666
$ws.Range("G20").Value = 13.59246626120456
$ws.Range("H20").Value = 9.572883886832349
$ws.Range("M20").Value = 68.898421753135
$ws.Range("O20").Value = 12.50212647198114

$ws.Range("C21").Value = 5.743232368273839
$ws.Range("D21").Value = 4.284492031256903
$ws.Range("E21").Value = 51.85224102274366
$ws.Range("F21").Value = 14.61437304988886
$ws.Range("G21").Value = 13.80067803502692
$ws.Range("H21").Value = 9.414497989792556
$ws.Range("M21").Value = 72.87164397591621
$ws.Range("O21").Value = 12.28660332758053

$ws.Range("C22").Value = 5.887144185041046
$ws.Range("D22").Value = 4.372323609465719
$ws.Range("E22").Value = 53.7375241598664
$ws.Range("F22").Value = 14.60089090269635
$ws.Range("G22").Value = 13.97589584786442
$ws.Range("H22").Value = 9.315975208509265
$ws.Range("M22").Value = 75.35922354516654
$ws.Range("O22").Value = 12.16008795330636

$ws.Range("C23").Value = 5.810657568247867
$ws.Range("D23").Value = 4.325688718631306
$ws.Range("E23").Value = 52.73880456633145
$ws.Range("F23").Value = 14.60626254466199
$ws.Range("G23").Value = 13.87877200103085
$ws.Range("H23").Value = 9.368093011105342
$ws.Range("M23").Value = 74.04148775961271
$ws.Range("O23").Value = 12.22625011161201

$ws.Range("C24").Value = 5.514144152954719
$ws.Range("D24").Value = 4.143803477352613
$ws.Range("E24").Value = 48.79065345282251
$ws.Range("F24").Value = 14.66586083025775
$ws.Range("G24").Value = 13.58971586624328
$ws.Range("H24").Value = 9.575554965997005
$ws.Range("M24").Value = 68.83162700858357
$ws.Range("O24").Value = 12.50588046049002

$ws.Range("C25").Value = 5.181391830316413
$ws.Range("D25").Value = 3.936836849150694
$ws.Range("E25").Value = 44.17132509431259
$ws.Range("F25").Value = 14.80860582717442
$ws.Range("G25").Value = 13.44108990975413
$ws.Range("H25").Value = 9.819930891755057
$ws.Range("M25").Value = 62.73827142210031
$ws.Range("O25").Value = 12.86435595207713
